# "add status in rm and pt"
# - pt (column C) goes from 2 -> 3 for the existing rows
# - rm_sources batch code (column H, shared string "0911-2SIDID") rolls from
#   09112SIDID -> 09113SIDID
# - a new day's worth of GI rows (6-9) is appended, mirroring rows 2-5 but
#   dated the next day (44056) with pt=3 and the new batch code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update pt (column C) for the existing rows from 2 to 3 ---
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 3

# --- append the new rows (6-9), mirroring rows 2-5 but for the next day ---
$ws.Range("A6").Value = 6800083954
$ws.Range("B6").Value = 44056
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "ASTB20_10"
$ws.Range("E6").Value = "TA57301F3C"
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = "SERAM"
$ws.Range("H6").Value = "09113SIDID"

$ws.Range("A7").Value = 6800083955
$ws.Range("B7").Value = 44056
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "HM20_B15"
$ws.Range("E7").Value = "TA57201F3C"
$ws.Range("F7").Value = 390
$ws.Range("G7").Value = "MOROTAI"
$ws.Range("H7").Value = "09113SIDID"

$ws.Range("A8").Value = 6800083955
$ws.Range("B8").Value = 44056
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "HM20_B15"
$ws.Range("E8").Value = "TA57301F3C"
$ws.Range("F8").Value = 5580
$ws.Range("G8").Value = "MOROTAI"
$ws.Range("H8").Value = "09113SIDID"

$ws.Range("A9").Value = 6800083955
$ws.Range("B9").Value = 44056
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "HM20_B15"
$ws.Range("E9").Value = "TA57401F3C"
$ws.Range("F9").Value = 1230
$ws.Range("G9").Value = "MOROTAI"
$ws.Range("H9").Value = "09113SIDID"

# give the new date cells the same short-date format as the existing ones
# (numFmtId 14 / "m/d/yy" -> reuses style index 1, same as B2:B5)
$ws.Range("B6:B9").NumberFormat = "m/d/yy"

# --- roll the rm_sources batch code forward on the existing rows too, so
# the shared string itself updates in place (09112SIDID -> 09113SIDID)
# instead of leaving a stray, now-unused old string behind ---
$ws.Range("H2").Value = "09113SIDID"
$ws.Range("H3").Value = "09113SIDID"
$ws.Range("H4").Value = "09113SIDID"
$ws.Range("H5").Value = "09113SIDID"

# widen column E (bmi_code) to fit the wider content, dropping autosize/bestFit
$ws.Columns("E").ColumnWidth = 15.8

# leave the selection on the last edited cell, like the saved workbook does
$ws.Range("E9").Select() | Out-Null
